# feat: Add language switch
#
# Inserts a "navbar.settings" entry right after "navbar.about" (pushing
# navbar.feedback / navbar.log / channel.* / appMenu.* down by one row),
# then appends a new "settings.language.*" block (title / auto / zh_cn /
# en_us) and one trailing blank row to keep the sheet's blank-row padding.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift rows 10-20 down to 11-21 in a single block copy (values +
#     formatting both travel with Copy(Destination)), freeing up row 10 for
#     the new "navbar.settings" entry. ---
$ws.Range("A10:C20").Copy($ws.Range("A11:C21"))

$ws.Range("A10").Value = "navbar.settings"
$ws.Range("B10").Value = "设置"
$ws.Range("C10").Value = "Settings"

# --- New "settings.language.*" block, rows 22-25. Copy the populated-data
#     style ("s=3" on A/B/C) down from row 20 first so the new cells match
#     the rest of the table, then fill in the text. ---
$ws.Range("A20:C20").Copy($ws.Range("A22:C25"))

$ws.Range("A22").Value = "settings.language.title"
$ws.Range("B22").Value = "语言"
$ws.Range("C22").Value = "Language"

$ws.Range("A23").Value = "settings.language.auto"
$ws.Range("B23").Value = "跟随系统"
$ws.Range("C23").Value = "Auto"

$ws.Range("A24").Value = "settings.language.zh_cn"
$ws.Range("B24").Value = "简体中文"
$ws.Range("C24").Value = "简体中文"

$ws.Range("A25").Value = "settings.language.en_us"
$ws.Range("B25").Value = "English"
$ws.Range("C25").Value = "English"

# --- New trailing blank row 29, matching the style + row height of the
#     existing blank row 28. ---
$ws.Range("A28:E28").Copy($ws.Range("A29:E29"))
$ws.Range("A29:E29").RowHeight = 20.1
